# Add basic gridline options (X Axis / Y Axis min/max/step) and update example excel

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 2 (pushes existing data down to start at row 4)
$ws.Rows("2:3").Insert()

# New gridline summary rows and header cells, written in the same order the
# original author entered them so shared-string indices line up
$ws.Cells.Item(3, 1).Value = "X Axis"
$ws.Cells.Item(1, 15).Value = "min"
$ws.Cells.Item(1, 16).Value = "max"
$ws.Cells.Item(2, 1).Value = "Y Axis"
$ws.Cells.Item(1, 17).Value = "step"

$ws.Cells.Item(2, 15).Value = -36
$ws.Cells.Item(2, 16).Value = 144
$ws.Cells.Item(2, 17).Value = 36

$ws.Cells.Item(3, 15).Value = 0
$ws.Cells.Item(3, 16).Value = 180
$ws.Cells.Item(3, 17).Value = 36

$ws.Range("Q4").Select()
